$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of column J (k value) over the data rows
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Summary rows 14-17
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Style the first summary value cell: bold, size 12, vertically centered
$r = $ws.Range("B14")
$r.Font.Bold = $true
$r.Font.Size = 12
$r.VerticalAlignment = -4108  # xlCenter

# Propagate the same formatting to the remaining summary cells via copy/paste-format
$r.Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row heights for the summary rows
$ws.Rows.Item(14).RowHeight = 15.6
$ws.Rows.Item(15).RowHeight = 15.6
$ws.Rows.Item(16).RowHeight = 15.6
$ws.Rows.Item(17).RowHeight = 15.6

# Selection as in the edited file
[void]$ws.Range("A14:B17").Select()

# Page setup (portrait, A4)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$wb.Save()
